# Restore the real calibrated parameter values that had been overwritten by
# stray "0.2" (and one "3.4") debugging values pushed into column B, which in
# turn had shifted the real values into column C. Move each column C true
# value back into column B and remove column C so the parameter_values sheet
# only has a single value column again for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = 0.058
$ws.Range("C18").ClearContents()

$ws.Range("B32").Value = 0.005
$ws.Range("C32").ClearContents()

$ws.Range("B34").Value = 0.01
$ws.Range("C34").ClearContents()

$ws.Range("B39").Value = 0.012
$ws.Range("C39").ClearContents()

$ws.Range("B40").Value = 0.01
$ws.Range("C40").ClearContents()

$ws.Range("B46").Value = 0.02
$ws.Range("C46").ClearContents()

$ws.Range("B47").Value = 0.184
$ws.Range("C47").ClearContents()

$ws.Range("B48").Value = 0.184
$ws.Range("C48").ClearContents()

$ws.Range("B49").Value = 0.33
$ws.Range("C49").ClearContents()

$ws.Range("B50").Value = 0.345
$ws.Range("C50").ClearContents()

$ws.Range("B62").Value = 0.01
$ws.Range("C62").ClearContents()

$ws.Range("B63").Value = 0.03
$ws.Range("C63").ClearContents()

$ws.Range("B66").Value = 0.05
$ws.Range("C66").ClearContents()

$ws.Range("B73").Value = 0.1
$ws.Range("C73").ClearContents()

$ws.Range("B74").Value = 0.184
$ws.Range("C74").ClearContents()

$ws.Range("B75").Value = 0.33
$ws.Range("C75").ClearContents()

# Row 33 only had a stray "3.4" debugging value in B with no real value
# pushed into C (the row legitimately has no numeric parameter here) -
# just remove it.
$ws.Range("B33").ClearContents()

# Restore the viewport/selection that was active when the author saved the
# file after making these corrections.
$ws.Range("B73:B75").Select()
